$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C) for rows 2-11 from 2023-10-13 (45212)
# to 2023-10-22 (45221), keeping existing number formatting/style intact.
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = 45221
}
